# Scheduled runner update: refresh computed profit/price figures on the
# Leve profit tables across the various crafting-class sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 913937
$ws.Range("I86").Value = 1113255.9
$ws.Range("J86").Value = 17002
$ws.Range("K86").Value = 1113255.9
$ws.Range("L86").Value = 17002
$ws.Range("M86").Value = -1112132.9
$ws.Range("N86").Value = -19248
$ws.Range("H89").Value = 913937
$ws.Range("I89").Value = 1113255.9
$ws.Range("J89").Value = 17002
$ws.Range("K89").Value = 5566279.5
$ws.Range("L89").Value = 85010
$ws.Range("M89").Value = -5560663.5
$ws.Range("N89").Value = -96242
$ws.Range("H137").Value = 1193.3489
$ws.Range("I137").Value = 923.5454999999999
$ws.Range("J137").Value = 1476
$ws.Range("K137").Value = 2770.6365
$ws.Range("L137").Value = 4428
$ws.Range("M137").Value = -220.6364999999996
$ws.Range("N137").Value = -9528

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H23").Value = 61253.25
$ws.Range("J23").Value = 58669
$ws.Range("L23").Value = 58669
$ws.Range("N23").Value = -59187
$ws.Range("H32").Value = 25126.2
$ws.Range("I32").Value = 4929.1777
$ws.Range("K32").Value = 4929.1777
$ws.Range("M32").Value = -4642.1777
$ws.Range("H37").Value = 4750
$ws.Range("J37").Value = 4750
$ws.Range("L37").Value = 4750
$ws.Range("N37").Value = -5296
$ws.Range("H44").Value = 16435.572
$ws.Range("J44").Value = 16435.572
$ws.Range("L44").Value = 16435.572
$ws.Range("N44").Value = -17411.572
$ws.Range("H55").Value = 8331.666999999999
$ws.Range("J55").Value = 8331.666999999999
$ws.Range("L55").Value = 8331.666999999999
$ws.Range("N55").Value = -8961.666999999999
$ws.Range("H80").Value = 10427
$ws.Range("J80").Value = 10427
$ws.Range("L80").Value = 10427
$ws.Range("N80").Value = -12423
$ws.Range("H83").Value = 10427
$ws.Range("J83").Value = 10427
$ws.Range("L83").Value = 31281
$ws.Range("N83").Value = -41265
$ws.Range("H132").Value = 12303050
$ws.Range("I132").Value = 22290328
$ws.Range("J132").Value = 779267.3
$ws.Range("K132").Value = 66870984
$ws.Range("L132").Value = 2337801.9
$ws.Range("M132").Value = -66868454
$ws.Range("N132").Value = -2342861.9

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 46668.332
$ws.Range("J4").Value = 46668.332
$ws.Range("L4").Value = 46668.332
$ws.Range("N4").Value = -46892.332
$ws.Range("H10").Value = 37530.75
$ws.Range("J10").Value = 70008
$ws.Range("L10").Value = 70008
$ws.Range("N10").Value = -70286
$ws.Range("H22").Value = 3532.6667
$ws.Range("I22").Value = 3532.6667
$ws.Range("K22").Value = 3532.6667
$ws.Range("M22").Value = -3182.6667
$ws.Range("H31").Value = 4312.986
$ws.Range("I31").Value = 3440.7896
$ws.Range("J31").Value = 5348.7188
$ws.Range("K31").Value = 3440.7896
$ws.Range("L31").Value = 5348.7188
$ws.Range("M31").Value = -3145.7896
$ws.Range("N31").Value = -5938.7188
$ws.Range("H32").Value = 6275
$ws.Range("I32").Value = 6700
$ws.Range("J32").Value = 5000
$ws.Range("K32").Value = 6700
$ws.Range("L32").Value = 5000
$ws.Range("M32").Value = -6384
$ws.Range("N32").Value = -5632
$ws.Range("H33").Value = 45741
$ws.Range("J33").Value = 45741
$ws.Range("L33").Value = 45741
$ws.Range("N33").Value = -46499
$ws.Range("H34").Value = 4312.986
$ws.Range("I34").Value = 3440.7896
$ws.Range("J34").Value = 5348.7188
$ws.Range("K34").Value = 3440.7896
$ws.Range("L34").Value = 5348.7188
$ws.Range("M34").Value = -3238.7896
$ws.Range("N34").Value = -5752.7188
$ws.Range("H36").Value = 57631.8
$ws.Range("I36").Value = 74000
$ws.Range("J36").Value = 46719.668
$ws.Range("K36").Value = 74000
$ws.Range("L36").Value = 46719.668
$ws.Range("M36").Value = -73612
$ws.Range("N36").Value = -47495.668
$ws.Range("H38").Value = 5000
$ws.Range("J38").Value = 5000
$ws.Range("L38").Value = 5000
$ws.Range("N38").Value = -5754
$ws.Range("H39").Value = 7441.8335
$ws.Range("I39").Value = 3550.3333
$ws.Range("J39").Value = 11333.333
$ws.Range("K39").Value = 3550.3333
$ws.Range("L39").Value = 11333.333
$ws.Range("M39").Value = -3159.3333
$ws.Range("N39").Value = -12115.333
$ws.Range("H40").Value = 57631.8
$ws.Range("I40").Value = 74000
$ws.Range("J40").Value = 46719.668
$ws.Range("K40").Value = 74000
$ws.Range("L40").Value = 46719.668
$ws.Range("M40").Value = -73840
$ws.Range("N40").Value = -47039.668
$ws.Range("H44").Value = 36142.332
$ws.Range("I44").Value = 12064
$ws.Range("J44").Value = 37862.215
$ws.Range("K44").Value = 12064
$ws.Range("L44").Value = 37862.215
$ws.Range("M44").Value = -11622
$ws.Range("N44").Value = -38746.215
$ws.Range("H45").Value = 12518.5
$ws.Range("I45").Value = 10000
$ws.Range("J45").Value = 13358
$ws.Range("K45").Value = 10000
$ws.Range("L45").Value = 13358
$ws.Range("M45").Value = -9407
$ws.Range("N45").Value = -14544
$ws.Range("H46").Value = 5000
$ws.Range("J46").Value = 5000
$ws.Range("L46").Value = 5000
$ws.Range("N46").Value = -5422
$ws.Range("H49").Value = 7441.8335
$ws.Range("I49").Value = 3550.3333
$ws.Range("J49").Value = 11333.333
$ws.Range("K49").Value = 3550.3333
$ws.Range("L49").Value = 11333.333
$ws.Range("M49").Value = -3368.3333
$ws.Range("N49").Value = -11697.333
$ws.Range("H50").Value = 40900
$ws.Range("I50").Value = 51350
$ws.Range("K50").Value = 51350
$ws.Range("M50").Value = -50725
$ws.Range("H51").Value = 20099
$ws.Range("J51").Value = 20099
$ws.Range("L51").Value = 20099
$ws.Range("N51").Value = -21571
$ws.Range("H56").Value = 61472.715
$ws.Range("I56").Value = 25000
$ws.Range("J56").Value = 76061.8
$ws.Range("K56").Value = 25000
$ws.Range("L56").Value = 76061.8
$ws.Range("M56").Value = -24155
$ws.Range("N56").Value = -77751.8
$ws.Range("H57").Value = 56020.332
$ws.Range("J57").Value = 62036.6
$ws.Range("L57").Value = 62036.6
$ws.Range("N57").Value = -63156.6
$ws.Range("H61").Value = 20099
$ws.Range("J61").Value = 20099
$ws.Range("L61").Value = 20099
$ws.Range("N61").Value = -20795
$ws.Range("H80").Value = 25064
$ws.Range("J80").Value = 25064
$ws.Range("L80").Value = 25064
$ws.Range("N80").Value = -27310
$ws.Range("H83").Value = 25064
$ws.Range("J83").Value = 25064
$ws.Range("L83").Value = 75192
$ws.Range("N83").Value = -86424
$ws.Range("H99").Value = 4373.6807
$ws.Range("I99").Value = 4294.731
$ws.Range("J99").Value = 4471.4287
$ws.Range("K99").Value = 4294.731
$ws.Range("L99").Value = 4471.4287
$ws.Range("M99").Value = -2796.731
$ws.Range("N99").Value = -7467.4287
$ws.Range("H100").Value = 57853.332
$ws.Range("J100").Value = 57853.332
$ws.Range("L100").Value = 57853.332
$ws.Range("N100").Value = -60017.332
$ws.Range("H126").Value = 4373.6807
$ws.Range("I126").Value = 4294.731
$ws.Range("J126").Value = 4471.4287
$ws.Range("K126").Value = 12884.193
$ws.Range("L126").Value = 13414.2861
$ws.Range("M126").Value = -10414.193
$ws.Range("N126").Value = -18354.2861

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 552
$ws.Range("J122").Value = 868.125
$ws.Range("L122").Value = 7813.125
$ws.Range("N122").Value = -12713.125
$ws.Range("H131").Value = 901.2381
$ws.Range("J131").Value = 1257.8334
$ws.Range("L131").Value = 3773.5002
$ws.Range("N131").Value = -13853.5002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2025
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 70002
$ws.Range("J2").Value = 70002
$ws.Range("L2").Value = 70002
$ws.Range("N2").Value = -70226
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H33").Value = 19481.7
$ws.Range("J33").Value = 53272.332
$ws.Range("L33").Value = 53272.332
$ws.Range("N33").Value = -53852.332
$ws.Range("H38").Value = 50000
$ws.Range("J38").Value = 50000
$ws.Range("L38").Value = 50000
$ws.Range("N38").Value = -50820
$ws.Range("H46").Value = 985.7143
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 966.6667
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 966.6667
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -1342.6667
$ws.Range("H50").Value = 18233.334
$ws.Range("J50").Value = 18233.334
$ws.Range("L50").Value = 18233.334
$ws.Range("N50").Value = -19507.334
$ws.Range("H56").Value = 52011
$ws.Range("I56").Value = 0
$ws.Range("J56").Value = 52011
$ws.Range("K56").Value = 0
$ws.Range("L56").Value = 52011
$ws.Range("M56").ClearContents()
$ws.Range("N56").Value = -53393
$ws.Range("H57").Value = 53973
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 53973
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 53973
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -55105
$ws.Range("H58").Value = 24697.666
$ws.Range("I58").Value = 4093
$ws.Range("J58").Value = 35000
$ws.Range("K58").Value = 4093
$ws.Range("L58").Value = 35000
$ws.Range("M58").Value = -3833
$ws.Range("N58").Value = -35520
$ws.Range("H136").Value = 401530.12
$ws.Range("I136").Value = 1251037.2
$ws.Range("J136").Value = 1762.0588
$ws.Range("K136").Value = 3753111.6
$ws.Range("L136").Value = 5286.1764
$ws.Range("M136").Value = -3750561.6
$ws.Range("N136").Value = -10386.1764

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 920.4
$ws.Range("I126").Value = 749.8333
$ws.Range("J126").Value = 1176.25
$ws.Range("K126").Value = 2249.4999
$ws.Range("L126").Value = 3528.75
$ws.Range("M126").Value = 220.5001000000002
$ws.Range("N126").Value = -8468.75
